$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Annual Population Survey rows (2-4): switch to latest "apsnew" live link,
#     and roll forward the release periods to the newest quarter ---
$ws.Range("B2").Value = "<a href='https://www.nomisweb.co.uk/datasets/apsnew'>Annual Population Survey</a>"
$ws.Range("C2").Value = "Apr 2022 - Mar 2023 (15/08/23)"
$ws.Range("D2").Value = "Jul 2022 - Jun 2023 (17/10/23)"

$ws.Range("C3").Value = "Apr 2022 - Mar 2023 (15/08/23)"
$ws.Range("D3").Value = "Jul 2022 - Jun 2023 (17/10/23)"

$ws.Range("B4").Value = "<a href='https://www.nomisweb.co.uk/datasets/apsnew'>Annual Population Survey</a>"
$ws.Range("C4").Value = "Apr 2022 - Mar 2023 (15/08/23)"
$ws.Range("D4").Value = "Jul 2022 - Jun 2023 (17/10/23)"

# --- Narrow column A, it no longer needs to fit the old (now removed) long text ---
$ws.Columns.Item(1).ColumnWidth = 85.83

# --- Update the selected cell to B5 ---
$ws.Range("B5").Select() | Out-Null
